$wb = $excel.ActiveWorkbook

# Insert "deleteCustomers" right after "test_suite"
$deleteCustomers = $wb.Worksheets.Add($null, $wb.Worksheets.Item("test_suite"))
$deleteCustomers.Name = "deleteCustomers"
$deleteCustomers.Range("A1").Value = "Name"
$deleteCustomers.Range("A2").Value = "Himun"
$deleteCustomers.Range("F9").Select() | Out-Null

# Append "customerLogin" as the last sheet (becomes the active tab)
$customerLogin = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$customerLogin.Name = "customerLogin"
$customerLogin.Range("A1").Value = "Name"
$customerLogin.Range("A2").Value = "Himun Trehan"
$customerLogin.Range("B2").Select() | Out-Null
